# Update database values and apply the new read_price algorithm results.
# This mirrors the source workbook's quarterly shift: every existing period
# column (D..M) shifts one column to the left (oldest quarter D dropped) and
# a brand-new trailing quarter is appended in column M. A handful of values
# that the refreshed pricing algorithm recomputed differ from a pure shift,
# so every cell below is written with its final, explicit value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: period headers (دوره مالی)
$row8 = @(
    "12 ماهه منتهی به 1399/09",
    "3 ماهه منتهی به 1399/12",
    "6 ماهه منتهی به 1400/03",
    "9 ماهه منتهی به 1400/06",
    "12 ماهه منتهی به 1400/09",
    "3 ماهه منتهی به 1400/12",
    "6 ماهه منتهی به 1401/03",
    "9 ماهه منتهی به 1401/06",
    "12 ماهه منتهی به 1401/09",
    "3 ماهه منتهی به 1401/12"
)

# Row 9: publish dates (تاریخ انتشار)
$row9 = @(
    "1400-12-04 (9)",
    "1401-01-30 (3)",
    "1401-05-30 (4)",
    "1401-07-27 (2)",
    "1402-01-16 (8)",
    "1402-01-30 (2)",
    "1401-05-30 (2)",
    "1401-07-27",
    "1402-01-30 (3)",
    "1402-01-30"
)

$rowData = @{
    8  = $row8
    9  = $row9
    11 = @(72582, 21516, 39578, 66533, 97818, 15819, 31572, 54386, 72542, 13988)
    12 = @(-54095, -18757, -35069, -57532, -85487, -15159, -31176, -52608, -69202, -12904)
    13 = @(18487, 2759, 4509, 9001, 12331, 660, 396, 1778, 3339, 1084)
    14 = @(-4082, -1164, -2387, -4646, -5823, -1226, -2543, -3708, -5141, -1097)
    15 = @("-", "-", "-", "-", "-", "-", "-", "-", "-", "-")
    16 = @(3389, 648, 1013, 1780, 2299, -336, -859, -321, -677, 145)
    17 = @(17794, 2243, 3135, 6135, 8807, -902, -3006, -2251, -2478, 132)
    18 = @(-2350, -927, -1996, -2823, -4164, -1163, -2488, -3394, -4552, -635)
    19 = @(131, 581, 78, 112, 221, 41, 117, 177, 248, 20)
    20 = @(15575, 1897, 1216, 3424, 4864, -2025, -5377, -5468, -6783, -483)
    21 = @(-2130, -347, 698, 314, -171, "-", -301, -293, -279, "-")
    22 = @(13445, 1551, 1914, 3738, 4693, -2025, -5678, -5760, -7062, -483)
    23 = @("-", "-", "-", "-", "-", "-", "-", "-", "-", "-")
    24 = @(13445, 1551, 1914, 3738, 4693, -2025, -5678, -5760, -7062, -483)
    25 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    26 = @(5421, 4422, 38613, 37403, 36000, 33479, 32444, 31517, 30042, 20288)
    27 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        # Column D is index 4 (1-based); the arrays above run D..M (10 columns).
        $ws.Cells.Item([int]$r, 4 + $i).Value = $vals[$i]
    }
}
